$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revisions to existing AgTests (F) / AgPosit (G) values for rows 618-684
$ws.Range("F618").Value = 38097
$ws.Range("G618").Value = 2677
$ws.Range("F643").Value = 43476
$ws.Range("F657").Value = 34066
$ws.Range("G657").Value = 874
$ws.Range("F658").Value = 27182
$ws.Range("F659").Value = 26213
$ws.Range("G659").Value = 849
$ws.Range("F660").Value = 6184
$ws.Range("F661").Value = 4876
$ws.Range("G661").Value = 314
$ws.Range("F662").Value = 12752
$ws.Range("F663").Value = 37087
$ws.Range("G663").Value = 1156
$ws.Range("F664").Value = 26416
$ws.Range("G664").Value = 776
$ws.Range("F665").Value = 28225
$ws.Range("G665").Value = 649
$ws.Range("F666").Value = 23757
$ws.Range("F667").Value = 16808
$ws.Range("G667").Value = 599
$ws.Range("F668").Value = 3376
$ws.Range("F669").Value = 23399
$ws.Range("F670").Value = 52405
$ws.Range("F671").Value = 32611
$ws.Range("G671").Value = 616
$ws.Range("F672").Value = 29744
$ws.Range("F673").Value = 10097
$ws.Range("F674").Value = 28631
$ws.Range("G674").Value = 681
$ws.Range("F675").Value = 13443
$ws.Range("F676").Value = 28079
$ws.Range("G676").Value = 449
$ws.Range("F677").Value = 55923
$ws.Range("F678").Value = 33809
$ws.Range("G678").Value = 524
$ws.Range("F679").Value = 29353
$ws.Range("G679").Value = 519
$ws.Range("F680").Value = 28108
$ws.Range("G680").Value = 547
$ws.Range("F681").Value = 26032
$ws.Range("G681").Value = 574
$ws.Range("F682").Value = 12359
$ws.Range("G682").Value = 408
$ws.Range("F683").Value = 23768
$ws.Range("G683").Value = 682
$ws.Range("F684").Value = 54829
$ws.Range("G684").Value = 1183

# Fill in the previously-missing AgTests/AgPosit for row 685
$ws.Range("F685").Value = 32436
$ws.Range("G685").Value = 1002

# New daily row 686
$ws.Range("A686").Value = 44580
$ws.Range("B686").Value = 896798
$ws.Range("C686").Value = 17803
$ws.Range("D686").Value = 6011
$ws.Range("E686").Value = 17520
$ws.Range("F686").Value = 23688
$ws.Range("G686").Value = 755
